$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add two new rows (12 and 13) to the table, duplicating row 11's
# --- formatting (they share the same style pattern per the target data).
$ws.Rows(11).Copy()
$ws.Rows(12).Insert()
$ws.Rows(11).Copy()
$ws.Rows(13).Insert()

# --- Populate the new rows' values.
# Row 12: Interoperability Solution / asset_release as solution in new collection
$ws.Range("F12").Value = "eProcurement"
$ws.Range("A12").Value = "Interoperability Solution"
$ws.Range("B12").Value = 58694
$ws.Range("C12").Value = "asset_release as solution in new collection"
$ws.Range("D12").Value = "New collection"
$ws.Range("G12").Value = "Yes"
$ws.Range("H12").Value = "Yes"

# Row 13: Project / project_project as solution in new collection
$ws.Range("F13").Value = "eProcurement"
$ws.Range("A13").Value = "Project"
$ws.Range("B13").Value = 26863
$ws.Range("C13").Value = "project_project as solution in new collection"
$ws.Range("D13").Value = "New collection"
$ws.Range("G13").Value = "Yes"
$ws.Range("H13").Value = "Yes"

# --- Update the "Migration rules" comments (column D) for the existing rows
# --- that used to reference the old, removed explanations.
$ws.Range("D10").Value = "More than one Community or Repository"
$ws.Range("D11").Value = "More than one Community or Repository"

$ws.Range("D3").Value = "No Repository or Community"
$ws.Range("D4").Value = "No Repository or Community"
$ws.Range("D5").Value = "No Repository or Community"
$ws.Range("D6").Value = "No Repository or Community"

# --- Resize the table / autofilter to include the two new rows.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A2:AA13"))

# --- Update the workbook-level filter database defined name.
foreach ($n in $wb.Names) {
    if ($n.Name -eq "1. Content items!_FilterDatabase") {
        $n.RefersTo = "='1. Content items'!`$A`$2:`$V`$13"
    }
}
